$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear D3 (was "Y") - no longer used
$ws.Range("D3").ClearContents()

# E4 label changes from "Radiant Flux (Watts)" to lowercase "radiant flux (Watts)"
$ws.Range("E4").Value = "radiant flux (Watts)"

# D5 gains a new "y" axis marker (previously empty)
$ws.Range("D5").Value = "y"

# Update the saved selection to E4
$ws.Range("E4").Select()
